$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 296, shifting rows 296:384 down to 297:385.
$ws.Rows.Item(296).Insert()

# Populate every column of the newly inserted row 296 with the new data
# record (same market/category metadata as its neighbours, new price data).
$ws.Cells.Item(296, 1).Value = 9                                      # A296 - Mercado ID
$ws.Cells.Item(296, 2).Value = "Vega Central Mapocho de Santiago"      # B296 - Mercado
$ws.Cells.Item(296, 3).Value = "Metropolitana"                        # C296 - Region
$ws.Cells.Item(296, 4).Value = 44627                                  # D296 - Fecha
$ws.Cells.Item(296, 5).Value = 13                                     # E296 - Codreg
$ws.Cells.Item(296, 6).Value = 100112012                              # F296 - Categoria ID
$ws.Cells.Item(296, 7).Value = "Espinaca"                             # G296 - Categoria
$ws.Cells.Item(296, 8).Value = "Sin especificar"                      # H296 - Variedad
$ws.Cells.Item(296, 9).Value = "Primera"                              # I296 - Calidad
$ws.Cells.Item(296, 10).Value = 61                                    # J296 - Volumen
$ws.Cells.Item(296, 11).Value = 12000                                 # K296 - Precio minimo
$ws.Cells.Item(296, 12).Value = 14000                                 # L296 - Precio maximo
$ws.Cells.Item(296, 13).Value = 13016                                 # M296 - Precio promedio ponderado
$ws.Cells.Item(296, 14).Value = "$/cuna 10 kilos"                     # N296 - Unidad de comercializacion
$ws.Cells.Item(296, 15).Value = "Provincia de Chacabuco"              # O296 - Origen
$ws.Cells.Item(296, 16).Value = 1302                                  # P296 - Precio $/Kg
$ws.Cells.Item(296, 17).Value = 10                                    # Q296 - Kg o Unidades
$ws.Cells.Item(296, 18).Value = "Hortaliza"                           # R296 - Clasificacion
